$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PBL II")

# --- Row height tweaks -------------------------------------------------
# Row 17: 77.25 -> 64.5
$ws.Rows.Item(17).RowHeight = 64.5
# Row 26: 26.25 -> default (AutoFit clears the explicit height)
$ws.Rows.Item(26).AutoFit() | Out-Null
# Row 34: 27.75 -> 16.5
$ws.Rows.Item(34).RowHeight = 16.5

# --- Re-style a handful of "Doing" cells in column F as "To release" ---
# (copy format from an existing "To release" cell, then restore the text)
$ws.Range("F35").Copy() | Out-Null
$ws.Range("F36").PasteSpecial(-4122) | Out-Null
$ws.Range("F36").Value2 = "To release"

$ws.Range("F54").PasteSpecial(-4122) | Out-Null
$ws.Range("F54").Value2 = "To release"
$ws.Range("D54").Value2 = "Bart"

$ws.Range("F55").PasteSpecial(-4122) | Out-Null
$ws.Range("F55").Value2 = "To release"

$ws.Range("F56").PasteSpecial(-4122) | Out-Null
$ws.Range("F56").Value2 = "To release"

$ws.Range("F57").PasteSpecial(-4122) | Out-Null
$ws.Range("F57").Value2 = "To release"

# --- Insert two new feedback rows before the trailing blank row --------
$ws.Range("A60:A61").EntireRow.Insert(-4121) | Out-Null

# Give the two new rows the same per-column formatting as row 59
$ws.Range("A59:F59").Copy() | Out-Null
$ws.Range("A60:F61").PasteSpecial(-4122) | Out-Null

# Column F on the new rows should look like "To release" (same as above)
$ws.Range("F35").Copy() | Out-Null
$ws.Range("F60:F61").PasteSpecial(-4122) | Out-Null

$ws.Range("A60").Value2 = 60
$ws.Range("B60").Value2 = "bug jobs, contact en footer sectie breder gemaakt (cfr services en solutions"
$ws.Range("C60").Value2 = ""
$ws.Range("D60").Value2 = "Bart"
$ws.Range("E60").Value2 = "HIGH"
$ws.Range("F60").Value2 = "To release"

$ws.Range("A61").Value2 = 61
$ws.Range("B61").Value2 = "zwitserse vlag weg bij referentie K&G als voorbeeld + ander hoover effect (shrink)"
$ws.Range("C61").Value2 = ""
$ws.Range("D61").Value2 = "Bart"
$ws.Range("E61").Value2 = "HIGH"
$ws.Range("F61").Value2 = "To release"

# --- Extend the autofilter / filter-database range to the new bottom ---
$ws.AutoFilterMode = $false
$ws.Range("A1:F63").AutoFilter() | Out-Null
$n = $wb.Names.Item("PBL II!_FilterDatabase")
$n.RefersTo = "='PBL II'!`$A`$1:`$F`$63"

# --- Restore the original selection (now pointing at the edited block) -
$ws.Range("B7:F61").Select() | Out-Null
